$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.848.00'
$ws.Range('E2').Value = '  +2.29%  '
$ws.Range('D3').Value = '3.376.00'
$ws.Range('E3').Value = '  +0.74%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '593.71'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +6.41%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '186.84'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.85%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.595'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.47%  '
$ws.Range('E9').Value = '  +2.21%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.591'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.11%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '47.63'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.08%  '
$ws.Range('E12').Value = '  +2.67%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '3.919.12'
$ws.Range('E13').Value = '  +0.89%  '
$ws.Range('B14').Value = 'BitcoinCash'
$ws.Range('C14').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '641.04'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +7.82%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.65'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.98%  '
$ws.Range('D16').Value = '67.795.86'
$ws.Range('E16').Value = '  +2.10%  '
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.119'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.27%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.376.02'
$ws.Range('E18').Value = '  +0.48%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '18.09'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.77%  '
$ws.Range('E20').Value = '  +1.01%  '
$ws.Range('E21').Value = '  +1.53%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '18.00'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.17%  '
$ws.Range('E23').Value = '  +1.73%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '99.90'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.76%  '
$ws.Range('E25').Value = '  +2.10%  '
$ws.Range('E26').Value = '  +4.84%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.80'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.12%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '32.66'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +6.32%  '
$ws.Range('E29').Value = '  +2.17%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.96'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.42%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '612.32'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.66%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.85'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.33%  '
$ws.Range('D33').Value = '4.020.61'
$ws.Range('E33').Value = '  +7.04%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.13'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.50%  '
$ws.Range('E35').Value = '  +2.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '56.28'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.09%  '
$ws.Range('E38').Value = '  +5.89%  '
$ws.Range('E39').Value = '  +3.82%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '33.81'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.12%  '
$ws.Range('E41').Value = '  +2.16%  '
$ws.Range('D42').Value = '0.0₃0706'
$ws.Range('E42').Value = '  +0.61%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.42'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.79%  '
$ws.Range('E44').Value = '  +1.17%  '
$ws.Range('E45').Value = '  +1.82%  '
$ws.Range('E46').Value = '  +0.53%  '
$ws.Range('E47').Value = '  +1.37%  '
$ws.Range('E48').Value = '  +0.12%  '
$ws.Range('E49').Value = '  +11.20%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '128.32'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.48%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.78'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +5.17%  '
